# Add a new quarterly column (BB) to the yoy_rt_data series:
#  - BB1 gets the new period's date serial, formatted like the other header
#    date cells (copy BA1's format so it keeps the custom date number format,
#    bold font, borders, and alignment).
#  - BB3:BB21 carry forward the latest known YoY value from column BA (the
#    "flat-line to most recent data point" pattern already used across the
#    sheet when a new column is appended before new actuals exist).
#  - Rows 2 and 22 only have a value in column A (start/end anchor rows), so
#    they are intentionally left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: new date column BB1 ---------------------------------------
$ws.Range("BA1").Copy() | Out-Null
$ws.Range("BB1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats - match BA1's style
$excel.CutCopyMode = $false
$ws.Range("BB1").Value = 45986

# --- Data rows: carry BA's last value into the new BB column ---------------
for ($row = 3; $row -le 21; $row++) {
    $baValue = $ws.Cells.Item($row, 53).Value2   # column 53 = BA
    $ws.Cells.Item($row, 54).Value = $baValue    # column 54 = BB
}
